$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 3 and row 4
$ws.Range("B3").Value = 0.77868337379253394
$ws.Range("C3").Value = 0.0060508705329897026

$ws.Range("B4").Value = 0.77868384931815227
$ws.Range("C4").Value = 0.00000061067904405150102

# Add new row 5 - force text type for the date-like string so it
# doesn't get auto-converted to a date serial number, then reset
# the style back to Normal/default (matches neighbouring date cells).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2024-04-11"
$ws.Range("A5").Style = "Normal"

$ws.Range("B5").Value = 0.77787814345819684
$ws.Range("C5").Value = -0.001034702158855535

# Update selection to match target state
$ws.Range("D4:D5").Select() | Out-Null
